$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price (D) values stay as text, matching the source data
# Volume(1h) (E) values already contain "%" and spaces so they remain text naturally.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.973.36'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.421.99'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.87'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.99'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '8.04'
$ws.Range('E9').Value = '  +3.06%  '
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('E11').Value = '  +3.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.005.98'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.50'
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.422.70'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.992.24'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.50'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.95'
$ws.Range('E20').Value = '  -3.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.81'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.572'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.559.23'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('E27').Value = '  -2.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.68'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('E31').Value = '  -3.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.23'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.51'
$ws.Range('E35').Value = '  +3.43%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.96'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '169.04'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '31.10'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.456.01'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0786'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.64'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.780'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.547.34'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.83'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.49'
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('E51').Value = '  -6.90%  '
